$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("A1:D111")
$key = $ws.Range("B1:B111")
$rng.Sort($key, 1, $null, $null, 1)
